$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Storage EFs")
$ws1 = $wb.Worksheets.Item("Slurry & application")
$ws4.Range("D2").Copy()
$ws1.Range("F5").PasteSpecial(-4122)
